$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bound-field expressions in row 3
$ws.Range("D3").Value = ":staffSet[].staffName"
$ws.Range("E3").Value = ":staffMap[].time.YMD"

# Update the current selection to match the merged cell B4:F5
$ws.Range("B4:F5").Select()
